$d = $word.ActiveDocument

# The paragraph that needs the recurring-task addition highlighted yellow.
$target = "Create new tags for priority assigning and grouping similar tasks."
$wdYellow = 7

foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    if ($r.Text -like "*$target*") {
        # Setting the Font's highlight (rather than just the Range) ensures
        # both the run text and the paragraph mark pick up the yellow
        # highlight, matching Word's own behaviour when highlighting a
        # whole paragraph (list item) via the UI.
        $r.Font.HighlightColorIndex = $wdYellow
        break
    }
}
